$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(17.65998168759154, 9.454334676803104, 8.66778661278788, 35.88709300095054, 3.673096781676145, 10.58216666743187, 11.5152360268251, 17.09381115951702, 19.39270941728943, 27.32680961507912),
    @(17.3116332721032, 9.277892576325129, 8.674903163824824, 35.95165870457684, 3.675260152760615, 10.60172772103869, 11.5206587173265, 17.02679434356548, 19.45187409895716, 27.38980771525483),
    @(17.09708237541352, 9.166855104773273, 8.680071783347728, 35.99940429315065, 3.676660041285918, 10.61435669875794, 11.52522818617098, 16.98785494458413, 19.49004931935115, 27.43440126564977),
    @(17.00960350535917, 9.120963388829553, 8.68237943563321, 36.02089473045378, 3.677248560986378, 10.61965906778992, 11.52740288121625, 16.97255334396941, 19.50607192364543, 27.4540572841419),
    @(16.99507836475678, 9.113305355158504, 8.682774797068946, 36.02458591731803, 3.677347376312282, 10.62054895495566, 11.52778289066672, 16.97004705129611, 19.50876063997787, 27.45741067637355),
    @(17.09590263817558, 9.166238747700067, 8.680102089075929, 35.99968589141595, 3.676667905093685, 10.6144275762898, 11.52525624812702, 16.98764627393585, 19.49026351778397, 27.43466034992736),
    @(17.54007950474394, 9.394079329168882, 8.670074850693181, 35.90767169064299, 3.673827891995985, 10.58878324107536, 11.51684892690961, 17.07025201049425, 19.4127266163704, 27.34730253173826),
    @(18.40076839660962, 9.817934208970891, 8.656730350093667, 35.79165721251665, 3.668823915983057, 10.54338076769011, 11.51016361792547, 17.24926367641339, 19.27528252320482, 27.22302903083926),
    @(19.02030188479462, 10.11346402566032, 8.65074978087701, 35.74585581721592, 3.66548847818383, 10.51297284671684, 11.51117620558312, 17.39044910071992, 19.18312684463478, 27.1605512141965),
    @(19.29808022924049, 10.24409655268802, 8.648852924514356, 35.73360376099321, 3.664044371388607, 10.49977365325583, 11.51291135726142, 17.456615181445, 19.14310224321219, 27.13841401845449),
    @(19.40258596185146, 10.29298843439057, 8.648252506240848, 35.73019903712675, 3.663507993365039, 10.49486609390562, 11.51375056421518, 17.48193518542143, 19.12821759732875, 27.13093647190626),
    @(19.38011068930118, 10.28248473799322, 8.648376583646238, 35.7308773822579, 3.663623046968851, 10.4959189977192, 11.51356174380211, 17.47647054977386, 19.1314112002426, 27.13250660806125),
    @(19.30669221881522, 10.24813061464144, 8.648801168599926, 35.73329890282761, 3.66400003363177, 10.4993680901982, 11.51297675560879, 17.45869306072514, 19.14187223308627, 27.13778068385324),
    @(19.26162947466719, 10.22701191664346, 8.649076573273863, 35.73494297582517, 3.664232311311889, 10.50149255714812, 11.5126421191836, 17.44783781999347, 19.14831528346681, 27.14112915113431),
    @(19.002058186638, 10.10484781917437, 8.650890272136486, 35.74682930221626, 3.665584322464564, 10.51384815792606, 11.51108836057488, 17.3861627890032, 19.18578063814436, 27.16212453766116),
    @(18.84170880648147, 10.02890848238978, 8.65221349554804, 35.75632020239332, 3.666432448611389, 10.52158988550947, 11.51046082908562, 17.34881377189296, 19.20924965441022, 27.17661537497866),
    @(18.74910445982368, 9.984873585632004, 8.653052148014687, 35.7625870189021, 3.666927161422594, 10.52610237691429, 11.51021992498915, 17.32751512045131, 19.22292707328593, 27.18554155051044),
    @(18.71768886396117, 9.969903782382893, 8.653349441017024, 35.76484758330201, 3.667095848174545, 10.52764048726233, 11.51015900869909, 17.3203357218432, 19.2275887299907, 27.18866532232235),
    @(18.85881789406185, 10.03702947421861, 8.652064612102462, 35.75522625842692, 3.66634145110644, 10.52075959436725, 11.51051521567815, 17.35277075022939, 19.20673285445132, 27.17501157673),
    @(19.32827630734563, 10.25823709961398, 8.648673263139608, 35.73255412834023, 3.663889019666951, 10.49835255012092, 11.51314364684062, 17.46390768446587, 19.13879220583242, 27.13620697816519),
    @(19.63107178233846, 10.39944080515549, 8.647143635463785, 35.72493409463576, 3.662347239768653, 10.48423666667247, 11.51592257784751, 17.53807585842019, 19.09597285915859, 27.11612339081167),
    @(19.46986382914691, 10.32439488813001, 8.647897377347913, 35.72834245913175, 3.663164550013721, 10.49172236416234, 11.51434270888988, 17.4983554994476, 19.11868177617107, 27.1263590379843),
    @(18.85108416779894, 10.03335914393735, 8.65213167958531, 35.75571830627113, 3.666382568915343, 10.52113477686628, 11.51049025402813, 17.35098125784676, 19.20787012407318, 27.17573480001087),
    @(18.16972675122485, 9.705928380381144, 8.659666779775284, 35.81612612982521, 3.670117483347272, 10.55514330835676, 11.51092789230593, 17.19908536861062, 19.31090944806617, 27.25159593639244)
)

$cols = @(2, 3, 4, 6, 7, 10, 12, 13, 14, 15)  # B, C, D, F, G, J, L, M, N, O

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowData = $data[$i]
    for ($j = 0; $j -lt $cols.Length; $j++) {
        $ws.Cells.Item($rowNum, $cols[$j]).Value = $rowData[$j]
    }
}

Write-Host "Updated loading_percent values for case with 380 kV"
